# Insert a new "Starting point" row as the first data row (row 2),
# pushing all existing data rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2, shifting rows 2..190 down to 3..191.
$ws.Rows.Item(2).Insert()

# Copy the header row's column-A formatting (bold, centered, bordered)
# onto the new A2 cell, matching the look of every other row's index column.
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)

# Populate the new row's values.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Starting point"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0

# Update the visible selection to span the full (now larger) data range.
$ws.Range("A2:A191").Select()
